$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values for the symbol list refresh.
# Leading apostrophe forces these numeric-looking strings to remain
# stored as text, matching the source data's inlineStr representation.
$ws.Range("D2").Value  = "'281.59"
$ws.Range("D3").Value  = "'20.64"
$ws.Range("D4").Value  = "'6.219"
$ws.Range("D5").Value  = "'0.06170"
$ws.Range("D6").Value  = "'3.586"
$ws.Range("D7").Value  = "'6.568"
$ws.Range("D8").Value  = "'1.498"
$ws.Range("D9").Value  = "'0.8186"
$ws.Range("D10").Value = "'0.01385"
$ws.Range("D12").Value = "'0.08371"
$ws.Range("D13").Value = "'0.03528"
$ws.Range("D14").Value = "'0.03215"
$ws.Range("D15").Value = "'0.09146"
$ws.Range("D17").Value = "'0.001639"
$ws.Range("D18").Value = "'0.04698"
$ws.Range("D19").Value = "'0.006437"
$ws.Range("D20").Value = "'0.006172"
$ws.Range("D23").Value = "'3.786"
$ws.Range("D25").Value = "'0.3355"
$ws.Range("D26").Value = "'0.1251"
$ws.Range("D40").Value = "'0.04704"
$ws.Range("D41").Value = "'0.007196"
$ws.Range("D42").Value = "'0.1100"
$ws.Range("D43").Value = "'0.003494"
$ws.Range("D44").Value = "'0.01109"
$ws.Range("D45").Value = "'0.00006498"
$ws.Range("D48").Value = "'0.002838"
$ws.Range("D49").Value = "'0.00001903"
